$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the "Price" column cells that are being rewritten to stay as TEXT
# (several new values, e.g. "1.001", "0.7928", "5.594", look like plain
# numbers to Excel's type inference and would otherwise be auto-converted
# to numeric values, losing the original text formatting used throughout
# this column).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.888.64"
$ws.Range("E2").Value = "  -0.33%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.896.85"

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - XRP
$ws.Range("D5").Value = "0.7928"
$ws.Range("E5").Value = "  -4.55%  "

# Row 6 - BNB
$ws.Range("D6").Value = "243.91"
$ws.Range("E6").Value = "  +0.83%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3162"
$ws.Range("E8").Value = "  -3.43%  "

# Row 9 - Solana
$ws.Range("D9").Value = "25.40"
$ws.Range("E9").Value = "  -4.10%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.07236"
$ws.Range("E10").Value = "  +3.05%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.29%  "

# Row 12 - now Polkadot (was Polygon)
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "5.594"
$ws.Range("E12").Value = "  +6.68%  "

# Row 13 - now Polygon (was Polkadot)
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7672"
$ws.Range("E13").Value = "  +0.16%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.860.60"
$ws.Range("E14").Value = "  -1.92%  "

# Row 15 - Litecoin
$ws.Range("E15").Value = "  +0.55%  "

# Row 16 - Uniswap
$ws.Range("D16").Value = "6.190"
$ws.Range("E16").Value = "  +5.90%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "29.864.61"
$ws.Range("E17").Value = "  -0.39%  "

# Row 18 - Avalanche
$ws.Range("D18").Value = "13.95"
$ws.Range("E18").Value = "  -0.98%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "244.66"
$ws.Range("E19").Value = "  +0.45%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.000007801"
$ws.Range("E20").Value = "  +0.71%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "8.213"
$ws.Range("E21").Value = "  +18.23%  "

# Row 22 - now Dai (was WrappedliquidstakedEther2.0)
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23 - now WrappedliquidstakedEther2.0 (was Dai)
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.140.69"
$ws.Range("E23").Value = "  -0.32%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  +0.11%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "0.1674"

# Row 26 - Cosmos
$ws.Range("D26").Value = "9.446"
$ws.Range("E26").Value = "  +2.05%  "

# Row 27 - Monero
$ws.Range("D27").Value = "164.27"
$ws.Range("E27").Value = "  -0.63%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  -0.96%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "2.066"
$ws.Range("E29").Value = "  -1.24%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "1.401"
$ws.Range("E30").Value = "  +3.27%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +2.52%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "4.488"
$ws.Range("E32").Value = "  +4.99%  "

# Row 33 - now InternetComputer(DFINITY) (was Hedera)
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.095"
$ws.Range("E33").Value = "  +0.79%  "

# Row 34 - now Hedera (was InternetComputer(DFINITY))
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.05544"
$ws.Range("E34").Value = "  -5.45%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "1.282"
$ws.Range("E35").Value = "  +1.34%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.7405"
$ws.Range("E36").Value = "  +1.33%  "

# Row 37 - Frax
$ws.Range("D37").Value = "0.9948"
$ws.Range("E37").Value = "  -0.38%  "

# Row 38 - HuobiToken
$ws.Range("D38").Value = "2.628"
$ws.Range("E38").Value = "  -3.32%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "0.01930"
$ws.Range("E39").Value = "  +0.77%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "2.782"
$ws.Range("E40").Value = "  +0.23%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.157.10"
$ws.Range("E41").Value = "  +16.62%  "

# Row 42 - Aave
$ws.Range("D42").Value = "74.44"
$ws.Range("E42").Value = "  +2.76%  "

# Row 43 - TheSandbox
$ws.Range("D43").Value = "0.4420"
$ws.Range("E43").Value = "  -0.46%  "

# Row 44 - FraxShare
$ws.Range("D44").Value = "5.900"
$ws.Range("E44").Value = "  +0.73%  "

# Row 45 - TrustWalletToken
$ws.Range("D45").Value = "0.8533"
$ws.Range("E45").Value = "  -0.45%  "

# Row 46 - Quant
$ws.Range("D46").Value = "104.74"
$ws.Range("E46").Value = "  +2.81%  "

# Row 47 - PaxDollar
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  +0.01%  "

# Row 48 - RenderToken
$ws.Range("D48").Value = "1.881"
$ws.Range("E48").Value = "  -0.99%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "10.02"
$ws.Range("E49").Value = "  +2.47%  "

# Row 50 - SynthetixNetwork
$ws.Range("D50").Value = "3.043"
$ws.Range("E50").Value = "  +11.73%  "

# Row 51 - Aptos
$ws.Range("D51").Value = "7.460"
$ws.Range("E51").Value = "  -1.14%  "
